# Add a new column 'Event' before the existing 'Event ' column on sheet Card18.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card18")

# Insert a new column at M, shifting existing M (and onward) to the right.
$ws.Columns("M").Insert()

# New header cell M1 (formatting is inherited automatically from the Insert).
$ws.Cells.Item(1, 13).Value = "Event"

# Fill data rows 2-12 of new column M with "nan" string values.
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 13).Value = "nan"
}
